# Fruta / hortaliza, semanal
# Insert a new weekly record as row 14 (pushing existing rows 14-33 down to 15-34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 14..33 down to 15..34, inserting a fresh (blank) row 14.
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with the new weekly record.
$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C14").Value = "Los Lagos"
$ws.Range("D14").Value = 44848
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 100112013
$ws.Range("G14").Value = "Alcachofa"
$ws.Range("H14").Value = "Española"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 200
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("N14").Value = "`$/caja 30 unidades"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 333
$ws.Range("Q14").Value = 30
$ws.Range("R14").Value = "Hortaliza"
